# Scheduled runner update: refresh market-board derived figures
# (currentAveragePrice / LevePrice / LeveProfit columns, H:N) across the
# leve-profit sheets. Values are plain numeric refreshes pulled from an
# external pricing source; some cells that no longer carry a meaningful
# value are cleared outright (matching upstream behaviour of omitting
# the cell rather than writing a blank/ zero placeholder), and a couple
# of previously-empty cells now receive a computed value.

$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
# Row 64 - Forged from the Void / Void Glue
$ws.Range("H64").Value = 7000
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 7000
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 7000
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -7496
# Row 67 - Dodging the Draft (L) / Void Glue
$ws.Range("H67").Value = 7000
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 7000
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 7000
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -8716
# Row 92 - Whinier than the Sword / Enchanted Koppranickel Ink
$ws.Range("H92").Value = 269.5
$ws.Range("I92").Value = 269.5
$ws.Range("K92").Value = 269.5
$ws.Range("M92").Value = 978.5
# Row 134 - Binding Spells / Crocodileskin Index
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
# Row 138 - All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 2888.2083
$ws.Range("J138").Value = 3868.8667
$ws.Range("L138").Value = 11606.6001
$ws.Range("N138").Value = -21886.6001

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32 - Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 1685.38
$ws.Range("I32").Value = 1432.6875
$ws.Range("K32").Value = 1432.6875
$ws.Range("M32").Value = -1145.6875
$ws.Range("N32").ClearContents()
# Row 61 - Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 2631.3125
$ws.Range("I61").Value = 1360.1
$ws.Range("K61").Value = 1360.1
$ws.Range("M61").Value = -1148.1
$ws.Range("N61").ClearContents()
# Row 109 - A Head of Demand / Deepgold Helm of Fending
$ws.Range("H109").Value = 74896.164
$ws.Range("J109").Value = 74896.164
$ws.Range("L109").Value = 74896.164
$ws.Range("N109").Value = -77670.164
# Row 132 - Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 3115.3076
$ws.Range("I132").Value = 1499.8334
$ws.Range("K132").Value = 4499.5002
$ws.Range("M132").Value = -1969.5002
$ws.Range("N132").ClearContents()
# Row 136 - Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 2631.3125
$ws.Range("I136").Value = 1360.1
$ws.Range("K136").Value = 4080.3
$ws.Range("M136").Value = -1530.3
$ws.Range("N136").ClearContents()

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
# Row 20 - Smelt and Dealt / Iron Ingot
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
# Row 22 - Driving Up the Wall / Elm Lumber
$ws.Range("H22").Value = 200
$ws.Range("J22").Value = 200
$ws.Range("L22").Value = 200
$ws.Range("N22").Value = -900
# Row 31 - Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 2137.4375
$ws.Range("I31").Value = 1183.25
$ws.Range("K31").Value = 1183.25
$ws.Range("M31").Value = -888.25
$ws.Range("N31").ClearContents()
# Row 34 - Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 2137.4375
$ws.Range("I34").Value = 1183.25
$ws.Range("K34").Value = 1183.25
$ws.Range("M34").Value = -981.25
$ws.Range("N34").ClearContents()
# Row 132 - Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 2865.3333
$ws.Range("I132").Value = 2865.3333
$ws.Range("K132").Value = 8595.999899999999
$ws.Range("M132").Value = -6065.999899999999
# Row 134 - Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 573.2
$ws.Range("J134").Value = 684.5
$ws.Range("L134").Value = 2053.5
$ws.Range("N134").Value = -7123.5

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
# Row 87 - Soup That Eats Like a Knight / Clam Chowder
$ws.Range("H87").Value = 389.5
$ws.Range("I87").Value = 389.5
$ws.Range("K87").Value = 1168.5
$ws.Range("M87").Value = 79.5
# Row 90 - Like Ma Used to Make (L) / Clam Chowder
$ws.Range("H90").Value = 389.5
$ws.Range("I90").Value = 389.5
$ws.Range("K90").Value = 3505.5
$ws.Range("M90").Value = 2734.5
# Row 123 - Topping Up the Pot / Zurek
$ws.Range("H123").Value = 4814.8
$ws.Range("I123").Value = 2222
$ws.Range("K123").Value = 6666
$ws.Range("M123").Value = -4216

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70 - Sky Is the Limit / Mythrite Ingot
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
# Row 73 - Hulls of Broken Dreams (L) / Mythrite Ingot
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
# Row 80 - Needs More Prayerbell / Hardsilver Ingot
$ws.Range("H80").Value = 3349
$ws.Range("J80").Value = 3399
$ws.Range("L80").Value = 3399
$ws.Range("N80").Value = -5395
# Row 83 - With a Noise That Reaches Heaven (L) / Hardsilver Ingot
$ws.Range("H83").Value = 3349
$ws.Range("J83").Value = 3399
$ws.Range("L83").Value = 16995
$ws.Range("N83").Value = -26979
# Row 132 - On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 4814
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 4814
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 14442
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -19502

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7 - Tan Before the Ban / Leather
$ws.Range("H7").Value = 1818
$ws.Range("I7").Value = 1818
$ws.Range("K7").Value = 1818
$ws.Range("M7").Value = -1706
# Row 68 - You Could Say It's a Moving Target / Wyvern Leather
$ws.Range("H68").Value = 6923.5
$ws.Range("J68").Value = 10198.25
$ws.Range("L68").Value = 10198.25
$ws.Range("N68").Value = -11696.25
# Row 71 - They Call It Bloody Mary (L) / Wyvern Leather
$ws.Range("H71").Value = 6923.5
$ws.Range("J71").Value = 10198.25
$ws.Range("L71").Value = 50991.25
$ws.Range("N71").Value = -58479.25
# Row 82 - Trainin' the Neck / Dragon Leather
$ws.Range("H82").Value = 75361.75
$ws.Range("J82").Value = 110003
$ws.Range("L82").Value = 110003
$ws.Range("N82").Value = -110725
# Row 85 - Training Is Only Skintight (L) / Dragon Leather
$ws.Range("H85").Value = 75361.75
$ws.Range("J85").Value = 110003
$ws.Range("L85").Value = 110003
$ws.Range("N85").Value = -112499
# Row 126 - Battered Books / Saiga Leather
$ws.Range("H126").Value = 1818
$ws.Range("I126").Value = 1818
$ws.Range("K126").Value = 5454
$ws.Range("M126").Value = -2984
# Row 132 - Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 3231.6667
$ws.Range("I132").Value = 1847.5
$ws.Range("K132").Value = 5542.5
$ws.Range("M132").Value = -3012.5
$ws.Range("N132").ClearContents()

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
# Row 74 - Clothing the Naked Truth / Ramie Robe of Casting
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
# Row 77 - When in Robes (L) / Ramie Robe of Casting
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
# Row 100 - Of Great Import / Kudzu Thread
$ws.Range("H100").Value = 1649.6666
$ws.Range("J100").Value = 1649.6666
$ws.Range("L100").Value = 3299.3332
$ws.Range("N100").Value = -4381.3332
# Row 113 - A Tender Table / Pixie Floss
$ws.Range("H113").Value = 882.1
$ws.Range("I113").Value = 796.5
$ws.Range("J113").Value = 1224.5
$ws.Range("K113").Value = 2389.5
$ws.Range("L113").Value = 3673.5
$ws.Range("M113").Value = -219.5
$ws.Range("N113").Value = -8013.5
# Row 132 - Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 3673.2104
$ws.Range("I132").Value = 1929.9
$ws.Range("K132").Value = 5789.700000000001
$ws.Range("M132").Value = -3259.700000000001
# Row 136 - Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 1541.8462
$ws.Range("I136").Value = 1614
$ws.Range("K136").Value = 4842
$ws.Range("M136").Value = -2292
